# Apply weekly CompStat data refresh to the CS-028 Precinct workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Header text updates (report volume/issue number and date range)
# ---------------------------------------------------------------------
$ws.Range("A8").Value = "Volume 32   Number  21"
$ws.Range("C9").Value = "Report Covering the Week  5/19/2025  Through  5/25/2025"

# ---------------------------------------------------------------------
# Row 14 - Murder
# ---------------------------------------------------------------------
$ws.Range("I14").Value = 3
$ws.Range("K14").Value = 50
$ws.Range("L14").Value = 0
$ws.Range("M14").Value = 50
$ws.Range("N14").Value = -76.923076923076

# ---------------------------------------------------------------------
# Row 15 - Rape
# ---------------------------------------------------------------------
$ws.Range("M15").Value = 0

# ---------------------------------------------------------------------
# Row 16 - Robbery
# ---------------------------------------------------------------------
$ws.Range("C16").Value = 3
$ws.Range("D16").Value = 5
$ws.Range("E16").Value = -40
$ws.Range("F16").Value = 9
$ws.Range("G16").Value = 12
$ws.Range("H16").Value = -25
$ws.Range("I16").Value = 50
$ws.Range("J16").Value = 75
$ws.Range("K16").Value = -33.333333333333
$ws.Range("L16").Value = -25.373134328358
$ws.Range("M16").Value = -38.271604938271
$ws.Range("N16").Value = -84.025559105431

# ---------------------------------------------------------------------
# Row 17 - Fel. Assault
# ---------------------------------------------------------------------
$ws.Range("D17").Value = 5
$ws.Range("E17").Value = 0
$ws.Range("F17").Value = 23
$ws.Range("G17").Value = 22
$ws.Range("H17").Value = 4.545454545454
$ws.Range("I17").Value = 93
$ws.Range("J17").Value = 104
$ws.Range("K17").Value = -10.576923076923
$ws.Range("L17").Value = -7.920792079207
$ws.Range("M17").Value = 24
$ws.Range("N17").Value = -63.671875

# ---------------------------------------------------------------------
# Row 18 - Burglary
# ---------------------------------------------------------------------
$ws.Range("D18").Value = 2
$ws.Range("F18").Value = 1
$ws.Range("G18").Value = 8
$ws.Range("H18").Value = -87.5
$ws.Range("J18").Value = 41
$ws.Range("K18").Value = -39.024390243902
$ws.Range("L18").Value = -48.979591836734
$ws.Range("M18").Value = -50
$ws.Range("N18").Value = -92.816091954023

# ---------------------------------------------------------------------
# Row 19 - Gr. Larceny
# ---------------------------------------------------------------------
$ws.Range("C19").Value = 8
$ws.Range("D19").Value = 11
$ws.Range("E19").Value = -27.272727272727
$ws.Range("F19").Value = 26
$ws.Range("G19").Value = 52
$ws.Range("H19").Value = -50
$ws.Range("I19").Value = 120
$ws.Range("J19").Value = 215
$ws.Range("K19").Value = -44.186046511627
$ws.Range("L19").Value = -13.043478260869
$ws.Range("M19").Value = 6.194690265486
$ws.Range("N19").Value = -25

# ---------------------------------------------------------------------
# Row 20 - G.L.A.  (C20 switches from a numeric 2 to the text "0")
# ---------------------------------------------------------------------
$ws.Cells.Item(14, 3).Copy()
$ws.Cells.Item(20, 3).PasteSpecial(-4122)
$ws.Cells.Item(14, 3).Copy()
$ws.Cells.Item(20, 3).PasteSpecial(-4163)
$ws.Application.CutCopyMode = $false

$ws.Range("D20").Value = 1
$ws.Range("E20").Value = -100
$ws.Range("G20").Value = 6
$ws.Range("H20").Value = -50
$ws.Range("J20").Value = 21
$ws.Range("K20").Value = -38.095238095238
$ws.Range("L20").Value = -40.909090909090
$ws.Range("M20").Value = 62.5
$ws.Range("N20").Value = -78.688524590163

# ---------------------------------------------------------------------
# Row 21 - TOTAL
# ---------------------------------------------------------------------
$ws.Range("C21").Value = 16
$ws.Range("D21").Value = 24
$ws.Range("E21").Value = -33.333333333333
$ws.Range("F21").Value = 63
$ws.Range("G21").Value = 102
$ws.Range("H21").Value = -38.235294117647
$ws.Range("I21").Value = 308
$ws.Range("J21").Value = 463
$ws.Range("K21").Value = -33.477321814254
$ws.Range("L21").Value = -19.371727748691
$ws.Range("M21").Value = -7.507507507507
$ws.Range("N21").Value = -73.516766981943

# ---------------------------------------------------------------------
# Row 22 - Transit
# ---------------------------------------------------------------------
$ws.Range("C22").Value = 1
$ws.Range("I22").Value = 9
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = 12.5

# ---------------------------------------------------------------------
# Row 23 - Housing
# ---------------------------------------------------------------------
$ws.Range("D23").Value = 3
$ws.Range("E23").Value = -66.666666666666
$ws.Range("F23").Value = 3
$ws.Range("G23").Value = 13
$ws.Range("H23").Value = -76.923076923076
$ws.Range("I23").Value = 27
$ws.Range("J23").Value = 42
$ws.Range("K23").Value = -35.714285714285
$ws.Range("L23").Value = -18.181818181818
$ws.Range("M23").Value = 107.692307692308

# ---------------------------------------------------------------------
# Row 24 - Petit Larceny
# ---------------------------------------------------------------------
$ws.Range("C24").Value = 32
$ws.Range("D24").Value = 22
$ws.Range("E24").Value = 45.454545454545
$ws.Range("F24").Value = 126
$ws.Range("G24").Value = 93
$ws.Range("H24").Value = 35.483870967741
$ws.Range("I24").Value = 635
$ws.Range("J24").Value = 491
$ws.Range("K24").Value = 29.327902240325
$ws.Range("L24").Value = 27.510040160642
$ws.Range("M24").Value = 71.621621621621

# ---------------------------------------------------------------------
# Row 25 - Retail Theft
# ---------------------------------------------------------------------
$ws.Range("C25").Value = 22
$ws.Range("D25").Value = 14
$ws.Range("E25").Value = 57.142857142857
$ws.Range("F25").Value = 91
$ws.Range("G25").Value = 48
$ws.Range("H25").Value = 89.583333333333
$ws.Range("I25").Value = 442
$ws.Range("J25").Value = 307
$ws.Range("K25").Value = 43.973941368078
$ws.Range("L25").Value = 45.394736842105

# ---------------------------------------------------------------------
# Row 26 - Misd. Assault
# ---------------------------------------------------------------------
$ws.Range("C26").Value = 9
$ws.Range("E26").Value = -25
$ws.Range("F26").Value = 43
$ws.Range("G26").Value = 36
$ws.Range("H26").Value = 19.444444444444
$ws.Range("I26").Value = 184
$ws.Range("J26").Value = 156
$ws.Range("K26").Value = 17.948717948717
$ws.Range("L26").Value = 9.523809523809
$ws.Range("M26").Value = -3.664921465968

# ---------------------------------------------------------------------
# Row 27 - UCR Rape*
# ---------------------------------------------------------------------
$ws.Range("I27").Value = 6
$ws.Range("K27").Value = -14.285714285714
$ws.Range("L27").Value = 20

# ---------------------------------------------------------------------
# Row 28 - Other Sex Crimes
# ---------------------------------------------------------------------
$ws.Range("F28").Value = 5
$ws.Range("G28").Value = 2
$ws.Range("H28").Value = 150
$ws.Range("L28").Value = 5.882352941176

# ---------------------------------------------------------------------
# Row 29 - Shooting Vic.
# ---------------------------------------------------------------------
$ws.Range("F29").Value = 2
$ws.Range("H29").Value = 100
$ws.Range("I29").Value = 5
$ws.Range("K29").Value = -28.571428571428
$ws.Range("L29").Value = -16.666666666666
$ws.Range("M29").Value = -44.444444444444
$ws.Range("N29").Value = -87.179487179487

# ---------------------------------------------------------------------
# Row 30 - Shooting Inc.
# ---------------------------------------------------------------------
$ws.Range("I30").Value = 4
$ws.Range("K30").Value = -33.333333333333
$ws.Range("L30").Value = -20
$ws.Range("M30").Value = -55.555555555555
$ws.Range("N30").Value = -88.235294117647
